$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append (dates, B, C, D, E, F, G)
$rows = @(
    @{ Date = "21-09-2021"; B = 10000; C = $null;  D = 0;    E = $null; F = $null; G = $null },
    @{ Date = "22-09-2021"; B = 10000; C = $null;  D = 0;    E = $null; F = $null; G = $null },
    @{ Date = "23-09-2021"; B = 10000; C = $null;  D = 0;    E = $null; F = $null; G = $null },
    @{ Date = "28-09-2021"; B = 10000; C = 15000;  D = 5000; E = 5000;  F = 0;     G = 3.72 },
    @{ Date = "30-09-2021"; B = 10000; C = $null;  D = 0;    E = $null; F = $null; G = $null }
)

$startRow = 25
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $data.Date
    $ws.Cells.Item($r, 2).Value = $data.B
    if ($null -ne $data.C) { $ws.Cells.Item($r, 3).Value = $data.C }
    $ws.Cells.Item($r, 4).Value = $data.D
    if ($null -ne $data.E) { $ws.Cells.Item($r, 5).Value = $data.E }
    if ($null -ne $data.F) { $ws.Cells.Item($r, 6).Value = $data.F }
    if ($null -ne $data.G) { $ws.Cells.Item($r, 7).Value = $data.G }
}
